$wb = $excel.ActiveWorkbook

# All four sheets share the same "Date" (K2) and "SequenceNo" (AG2) /
# "DateandTime" (N2) cells that were bumped forward to a later test run.
# K2's new text ("02-06-2024") parses as a valid MM-DD-YYYY date under
# en-US locale rules, so it must be forced to stay literal text (same as
# how the workbook already stores it) via a leading apostrophe - Excel
# strips the apostrophe itself and only uses it to flag quoted input.

$sheet1 = $wb.Sheets.Item(1)
$sheet2 = $wb.Sheets.Item(2)
$sheet3 = $wb.Sheets.Item(3)
$sheet4 = $wb.Sheets.Item(4)

# --- Sheet1 ("New Add" scenario) ---
$sheet1.Range("K2").Value = "'02-06-2024"
$sheet1.Range("N2").Value = "30-05-2024 06:07:11 PM"
$sheet1.Range("AG2").Value = "ET478"
$sheet1.Range("AK2").Value = "'3"

# --- Sheet2 ("Summary Add" scenario) ---
$sheet2.Range("K2").Value = "'02-06-2024"
$sheet2.Range("N2").Value = "30-05-2024 06:12:52 PM"
$sheet2.Range("AG2").Value = "ET476"

# --- Sheet3 ("Duplicate Add" scenario) ---
$sheet3.Range("K2").Value = "'02-06-2024"
$sheet3.Range("N2").Value = "30-05-2024 06:12:52 PM"
$sheet3.Range("AG2").Value = "ET477"

# --- Sheet4 ("Edit Save" scenario) ---
$sheet4.Range("K2").Value = "'02-06-2024"
$sheet4.Range("N2").Value = "30-05-2024 06:12:52 PM"
# AG2 on Sheet4 keeps the same displayed value (ET164); nothing to change.

# The active sheet's selection moved from A2 to D10 before the file was saved.
$sheet1.Range("D10").Select() | Out-Null
